$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Julio de 2020 a las 22:36"

# Update country stat rows (Casos totales, Nuevos casos, Casos activos, Recuperados, Muertes hoy, Muertes)
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4416832
$ws.Range("C4").Value = 44993
$ws.Range("D4").Value = 2121993
$ws.Range("E4").Value = 2144627
$ws.Range("G4").Value = 364
$ws.Range("H4").Value = 150212

# Row 6 - India
$ws.Range("B6").Value = 1482503
$ws.Range("C6").Value = 46484
$ws.Range("D6").Value = 953189
$ws.Range("E6").Value = 495866

# Row 8 - Sudafrica
$ws.Range("B8").Value = 452529
$ws.Range("C8").Value = 7096
$ws.Range("D8").Value = 274925
$ws.Range("E8").Value = 170537
$ws.Range("G8").Value = 298
$ws.Range("H8").Value = 7067

# Row 12 - España
$ws.Range("B12").Value = 325862
$ws.Range("C12").Value = 2120
$ws.Range("H12").Value = 28434

# Row 21 - Alemania
$ws.Range("B21").Value = 207372
$ws.Range("C21").Value = 631
$ws.Range("E21").Value = 7568
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 9204

# Row 22 - Francia
$ws.Range("B22").Value = 183079
$ws.Range("C22").Value = 514
$ws.Range("E22").Value = 72055
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = 30209

# Row 32 - Ecuador
$ws.Range("B32").Value = 81161
$ws.Range("C32").Value = 467
$ws.Range("E32").Value = 40733
$ws.Range("G32").Value = 17
$ws.Range("H32").Value = 5532

# Row 68 - Republica Dominicana
$ws.Range("B68").Value = 17975
$ws.Range("C68").Value = 372
$ws.Range("D68").Value = 7833
$ws.Range("E68").Value = 9857
$ws.Range("G68").Value = 5
$ws.Range("H68").Value = 285
